$d = $word.ActiveDocument

$d.Content.Find.Execute("12÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "17÷9=", 2) | Out-Null
$d.Content.Find.Execute("73÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "10÷9=", 2) | Out-Null
$d.Content.Find.Execute("79÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "18÷3=", 2) | Out-Null
$d.Content.Find.Execute("26÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "22÷9=", 2) | Out-Null
$d.Content.Find.Execute("47÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "16÷3=", 2) | Out-Null
$d.Content.Find.Execute("32÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "29÷3=", 2) | Out-Null
$d.Content.Find.Execute("69÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "96÷7=", 2) | Out-Null
$d.Content.Find.Execute("53÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "86÷2=", 2) | Out-Null
$d.Content.Find.Execute("45÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "55÷7=", 2) | Out-Null
$d.Content.Find.Execute("35÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "61÷6=", 2) | Out-Null
$d.Content.Find.Execute("52÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "45÷9=", 2) | Out-Null
$d.Content.Find.Execute("23÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "36÷8=", 2) | Out-Null
$d.Content.Find.Execute("32÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "47÷6=", 2) | Out-Null
$d.Content.Find.Execute("25÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "60÷2=", 2) | Out-Null
$d.Content.Find.Execute("37÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "81÷6=", 2) | Out-Null
$d.Content.Find.Execute("62÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "88÷5=", 2) | Out-Null
$d.Content.Find.Execute("63÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "62÷4=", 2) | Out-Null
$d.Content.Find.Execute("37÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "80÷6=", 2) | Out-Null
$d.Content.Find.Execute("50÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "36÷3=", 2) | Out-Null
$d.Content.Find.Execute("10÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "48÷9=", 2) | Out-Null
$d.Content.Find.Execute("41÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "41÷9=", 2) | Out-Null
$d.Content.Find.Execute("19÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "24÷7=", 2) | Out-Null
$d.Content.Find.Execute("40÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "21÷5=", 2) | Out-Null
$d.Content.Find.Execute("98÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "87÷7=", 2) | Out-Null
$d.Content.Find.Execute("77÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "48÷3=", 2) | Out-Null
